# Register send Verify email and Register User from Excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite header row ---
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Role"

# --- Row 2: Doru, registered via emailed link ---
$ws.Range("A2").Value = "bocaioandoru12@gmail.com"
$ws.Range("B2").Value = "Doru"
$ws.Range("C2").Value = "Student"

# --- Row 3: Alex ---
$ws.Range("A3").Value = "birlea24@gmail.com"
$ws.Range("B3").Value = "Alex"
$ws.Range("C3").Value = "Student"

# --- Drop the old extra rows (4-7) so the table is only 3 rows tall ---
$ws.Rows("4:7").Delete()

# --- Turn the e-mail in A2 into a live "verify e-mail" hyperlink ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:bocaioandoru12@gmail.com")

# --- Column sizing: widen A for the e-mail address, narrow C ---
$ws.Columns("A").ColumnWidth = 24.77734375
$ws.Columns("C").ColumnWidth = 8

# --- Move the active selection like the saved workbook had it ---
[void]$ws.Range("E3").Select()
